$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 2.14
$ws.Range("L3").Value = 1.01
$ws.Range("N3").Value = 2.78
$ws.Range("O3").Value = 1.32
$ws.Range("S3").Value = 3.35
$ws.Range("T3").Value = 1.83
$ws.Range("U3").Value = 1.82
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 980
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 980
$ws.Range("AD3").Value = 980
$ws.Range("AG3").Value = 980
$ws.Range("AI3").Value = 65
$ws.Range("AK3").Value = 70
$ws.Range("AO3").Value = 1000

$ws.Range("F4").Value = 2.16
$ws.Range("I4").Value = 3.75
$ws.Range("O4").Value = 1.27
$ws.Range("W4").Value = 1.84
$ws.Range("Z4").Value = 30

$ws.Range("J5").Value = 3.15
$ws.Range("V5").Value = 1.3
$ws.Range("AL5").Value = 980

$ws.Range("H7").Value = 2.34
$ws.Range("Q7").Value = 2.1

$ws.Range("L8").Value = 1.33

$ws.Range("F10").Value = 2.2
$ws.Range("J10").Value = 3.9

$ws.Range("F11").Value = 5.4
$ws.Range("G11").Value = 5.5
$ws.Range("H11").Value = 1.69
$ws.Range("I11").Value = 1.71
$ws.Range("O11").Value = 1.23
$ws.Range("P11").Value = 2.34
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 2.8
$ws.Range("V11").Value = 2.42
$ws.Range("W11").Value = 1.22
$ws.Range("Y11").Value = 10.5
$ws.Range("AA11").Value = 17
$ws.Range("AC11").Value = 9.6
$ws.Range("AG11").Value = 20
$ws.Range("AO11").Value = 8

$ws.Range("N12").Value = 3.85
$ws.Range("S12").Value = 3.55

$ws.Range("G13").Value = 1.89

$ws.Range("F14").Value = 1.76
$ws.Range("G14").Value = 1.77
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 5.1
$ws.Range("S14").Value = 2.96
$ws.Range("V14").Value = 1.24
$ws.Range("AL14").Value = 29

$ws.Range("H15").Value = 4.6
$ws.Range("V15").Value = 1.23
$ws.Range("W15").Value = 1.95
$ws.Range("AA15").Value = 150
$ws.Range("AE15").Value = 90
$ws.Range("AI15").Value = 100

$ws.Range("F16").Value = 1.25
$ws.Range("G16").Value = 980
$ws.Range("H16").Value = 1.04
$ws.Range("J16").Value = 1.25
$ws.Range("V16").Value = 1.01

$ws.Range("J17").Value = 3.6
$ws.Range("N17").Value = 3.75
$ws.Range("O17").Value = 1.3
$ws.Range("Q17").Value = 1.89
$ws.Range("R17").Value = 1.38
$ws.Range("U17").Value = 2.14
$ws.Range("W17").Value = 1.8
